# A new weekly price record was inserted in the middle of the "Espinaca"
# price series (at what is currently row 131), pushing every subsequent
# record down by one row. This grows the sheet from 189 to 190 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 131; Excel shifts rows 131-189 down
# to 132-190 automatically (and grows the sheet dimension accordingly).
$ws.Rows("131:131").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(131, 1).Value  = 3
$ws.Cells.Item(131, 2).Value  = 'Femacal de La Calera'
$ws.Cells.Item(131, 3).Value  = 'Coquimbo'
$ws.Cells.Item(131, 4).Value  = 44466
$ws.Cells.Item(131, 5).Value  = 5
$ws.Cells.Item(131, 6).Value  = 100112012
$ws.Cells.Item(131, 7).Value  = 'Espinaca'
$ws.Cells.Item(131, 8).Value  = 'Sin especificar'
$ws.Cells.Item(131, 9).Value  = 'Primera'
$ws.Cells.Item(131, 10).Value = 230
$ws.Cells.Item(131, 11).Value = 2300
$ws.Cells.Item(131, 12).Value = 2500
$ws.Cells.Item(131, 13).Value = 2396
$ws.Cells.Item(131, 14).Value = '$/docena de atados (3 kilos)'
$ws.Cells.Item(131, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(131, 16).Value = 799
$ws.Cells.Item(131, 17).Value = 3
$ws.Cells.Item(131, 18).Value = 'Hortaliza'
